$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new empty row at row 3, shifting existing rows 3..15 down to 4..16.
$ws.Rows("3:3").Insert()

# Reflect the active selection recorded in the saved file.
$ws.Range("B5").Select()
